# Update report header text (mayor, volume/number, reporting week dates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Value = "Volume 31   Number  39"
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# Helper: write a numeric value and pin its number format so the cell
# reuses the workbook's existing style (matches canonical output).
function Set-NumCell($ws, $addr, $val, $fmt) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $fmt
}

# Helper: write a literal text value (e.g. "0" placeholder or "***.*"
# suppressed-value marker) and restore formatting from a known text-style cell.
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

$fmtInt = "#,##0"
$fmtPct1 = "#,##0.0;""-""#,##0.0"
$fmtPct2 = "#,##0.00;""-""#,##0.00"

# --- Weekly / 28-day / YTD / 2-year crime count updates ---
Set-NumCell $ws "L15" 66.666666666666 $fmtPct1
Set-NumCell $ws "M15" 11.111111111111 $fmtPct1
Set-NumCell $ws "C16" 1 $fmtInt
Set-NumCell $ws "F16" 6 $fmtInt
Set-NumCell $ws "G16" 13 $fmtInt
Set-NumCell $ws "H16" -53.846153846153 $fmtPct1
Set-NumCell $ws "I16" 107 $fmtInt
Set-NumCell $ws "K16" 11.458333333333 $fmtPct1
Set-NumCell $ws "L16" -15.079365079365 $fmtPct1
Set-NumCell $ws "M16" -27.210884353741 $fmtPct1
Set-NumCell $ws "N16" -83.90977443609 $fmtPct1
Set-NumCell $ws "C17" 5 $fmtInt
Set-NumCell $ws "D17" 1 $fmtInt
Set-NumCell $ws "E17" 400 $fmtPct1
Set-NumCell $ws "F17" 16 $fmtInt
Set-NumCell $ws "G17" 10 $fmtInt
Set-NumCell $ws "H17" 60 $fmtPct1
Set-NumCell $ws "I17" 116 $fmtInt
Set-NumCell $ws "J17" 138 $fmtInt
Set-NumCell $ws "K17" -15.942028985507 $fmtPct1
Set-NumCell $ws "L17" -15.942028985507 $fmtPct1
Set-NumCell $ws "M17" 48.717948717948 $fmtPct1
Set-NumCell $ws "N17" -68.648648648648 $fmtPct1
Set-NumCell $ws "D18" 2 $fmtInt
Set-NumCell $ws "E18" -50 $fmtPct1
Set-NumCell $ws "G18" 8 $fmtInt
Set-NumCell $ws "H18" 0 $fmtPct1
Set-NumCell $ws "I18" 83 $fmtInt
Set-NumCell $ws "J18" 130 $fmtInt
Set-NumCell $ws "K18" -36.153846153846 $fmtPct1
Set-NumCell $ws "L18" -40.287769784172 $fmtPct1
Set-NumCell $ws "M18" -14.432989690721 $fmtPct1
Set-NumCell $ws "N18" -89.572864321608 $fmtPct1
Set-NumCell $ws "D19" 4 $fmtInt
Set-NumCell $ws "E19" 300 $fmtPct1
Set-NumCell $ws "F19" 45 $fmtInt
Set-NumCell $ws "G19" 30 $fmtInt
Set-NumCell $ws "H19" 50 $fmtPct1
Set-NumCell $ws "I19" 381 $fmtInt
Set-NumCell $ws "J19" 375 $fmtInt
Set-NumCell $ws "K19" 1.6 $fmtPct1
Set-NumCell $ws "L19" -4.271356783919 $fmtPct1
Set-NumCell $ws "M19" -2.056555269922 $fmtPct1
Set-NumCell $ws "N19" -52.670807453416 $fmtPct1
Set-NumCell $ws "D20" 3 $fmtInt
Set-NumCell $ws "E20" -33.333333333333 $fmtPct1
Set-NumCell $ws "F20" 7 $fmtInt
Set-NumCell $ws "G20" 9 $fmtInt
Set-NumCell $ws "H20" -22.222222222222 $fmtPct1
Set-NumCell $ws "I20" 52 $fmtInt
Set-NumCell $ws "J20" 79 $fmtInt
Set-NumCell $ws "K20" -34.177215189873 $fmtPct1
Set-NumCell $ws "L20" -26.760563380281 $fmtPct1
Set-NumCell $ws "M20" 79.310344827586 $fmtPct1
Set-NumCell $ws "N20" -93.020134228187 $fmtPct1
Set-NumCell $ws "C21" 25 $fmtInt
Set-NumCell $ws "D21" 10 $fmtInt
Set-NumCell $ws "E21" 150 $fmtPct2
Set-NumCell $ws "F21" 82 $fmtInt
Set-NumCell $ws "G21" 70 $fmtInt
Set-NumCell $ws "H21" 17.142857142857 $fmtPct2
Set-NumCell $ws "I21" 749 $fmtInt
Set-NumCell $ws "J21" 830 $fmtInt
Set-NumCell $ws "K21" -9.759036144578 $fmtPct2
Set-NumCell $ws "L21" -14.789533560864 $fmtPct2
Set-NumCell $ws "M21" -0.266311584553 $fmtPct2
Set-NumCell $ws "N21" -78.137769994162 $fmtPct2
Set-NumCell $ws "D22" 1 $fmtInt
Set-NumCell $ws "E22" 0 $fmtPct1
Set-NumCell $ws "F22" 3 $fmtInt
Set-NumCell $ws "G22" 2 $fmtInt
Set-NumCell $ws "H22" 50 $fmtPct1
Set-NumCell $ws "I22" 24 $fmtInt
Set-NumCell $ws "J22" 16 $fmtInt
Set-NumCell $ws "K22" 50 $fmtPct1
Set-NumCell $ws "L22" 20 $fmtPct1
Set-NumCell $ws "M22" 33.333333333333 $fmtPct1
Set-NumCell $ws "C23" 4 $fmtInt
Set-NumCell $ws "F23" 10 $fmtInt
Set-NumCell $ws "G23" 6 $fmtInt
Set-NumCell $ws "H23" 66.666666666666 $fmtPct1
Set-NumCell $ws "I23" 74 $fmtInt
Set-NumCell $ws "K23" -10.843373493975 $fmtPct1
Set-NumCell $ws "L23" -20.430107526881 $fmtPct1
Set-NumCell $ws "M23" 27.586206896551 $fmtPct1
Set-NumCell $ws "C24" 37 $fmtInt
Set-NumCell $ws "D24" 33 $fmtInt
Set-NumCell $ws "E24" 12.121212121212 $fmtPct1
Set-NumCell $ws "F24" 127 $fmtInt
Set-NumCell $ws "G24" 103 $fmtInt
Set-NumCell $ws "H24" 23.300970873786 $fmtPct1
Set-NumCell $ws "I24" 896 $fmtInt
Set-NumCell $ws "J24" 1202 $fmtInt
Set-NumCell $ws "K24" -25.457570715474 $fmtPct1
Set-NumCell $ws "L24" -35.400144196106 $fmtPct1
Set-NumCell $ws "M24" 11.581569115815 $fmtPct1
Set-NumCell $ws "C25" 25 $fmtInt
Set-NumCell $ws "D25" 18 $fmtInt
Set-NumCell $ws "E25" 38.888888888888 $fmtPct1
Set-NumCell $ws "F25" 73 $fmtInt
Set-NumCell $ws "G25" 71 $fmtInt
Set-NumCell $ws "H25" 2.81690140845 $fmtPct1
Set-NumCell $ws "I25" 494 $fmtInt
Set-NumCell $ws "J25" 823 $fmtInt
Set-NumCell $ws "K25" -39.975698663426 $fmtPct1
Set-NumCell $ws "L25" -44.556677890011 $fmtPct1
Set-NumCell $ws "C26" 5 $fmtInt
Set-NumCell $ws "D26" 2 $fmtInt
Set-NumCell $ws "E26" 150 $fmtPct1
Set-NumCell $ws "F26" 27 $fmtInt
Set-NumCell $ws "H26" 12.5 $fmtPct1
Set-NumCell $ws "I26" 237 $fmtInt
Set-NumCell $ws "J26" 222 $fmtInt
Set-NumCell $ws "K26" 6.756756756756 $fmtPct1
Set-NumCell $ws "L26" 2.597402597402 $fmtPct1
Set-NumCell $ws "M26" -8.13953488372 $fmtPct1
Set-NumCell $ws "C27" 1 $fmtInt
Set-NumCell $ws "D27" 1 $fmtInt
Set-NumCell $ws "E27" 0 $fmtPct1
Set-NumCell $ws "F27" 1 $fmtInt
Set-NumCell $ws "G27" 2 $fmtInt
Set-NumCell $ws "H27" -50 $fmtPct1
Set-NumCell $ws "I27" 17 $fmtInt
Set-NumCell $ws "J27" 16 $fmtInt
Set-NumCell $ws "K27" 6.25 $fmtPct1
Set-NumCell $ws "L27" 70 $fmtPct1
Set-NumCell $ws "L28" -21.621621621621 $fmtPct1
Set-NumCell $ws "F31" 1 $fmtInt
Set-NumCell $ws "H31" 0 $fmtPct1
Set-NumCell $ws "I31" 18 $fmtInt
Set-NumCell $ws "K31" 63.636363636363 $fmtPct1
Set-NumCell $ws "L31" 12.5 $fmtPct1

# --- Cells reverting to placeholder text ("0" counts / "***.*" suppressed %) ---
Set-TextCell $ws "D16" "0"
Set-TextCell $ws "E16" "***.*"
Set-TextCell $ws "D23" "0"
Set-TextCell $ws "E23" "***.*"
Set-TextCell $ws "C28" "0"
Set-TextCell $ws "D28" "0"
Set-TextCell $ws "E28" "***.*"
Set-TextCell $ws "D29" "0"
Set-TextCell $ws "E29" "***.*"
Set-TextCell $ws "D30" "0"
Set-TextCell $ws "E30" "***.*"
